$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.37287688724919
$ws.Range("C2").Value = 10.88001358099455
$ws.Range("D2").Value = 4.971457173895673
$ws.Range("F2").Value = 26.97504469970814
$ws.Range("G2").Value = 33.82363493010849
$ws.Range("H2").Value = 15.20035226987862
$ws.Range("I2").Value = 23.49204998432765
$ws.Range("L2").Value = 10.83551531785542
$ws.Range("M2").Value = 15.29865013439652
$ws.Range("N2").Value = 18.16891629254204

$ws.Range("B3").Value = 14.89161267742864
$ws.Range("C3").Value = 10.52710311849199
$ws.Range("D3").Value = 4.967978970695796
$ws.Range("F3").Value = 26.8725176172955
$ws.Range("G3").Value = 33.59216657130011
$ws.Range("H3").Value = 15.22497386548312
$ws.Range("I3").Value = 23.56375886033323
$ws.Range("L3").Value = 10.85095501459578
$ws.Range("M3").Value = 15.21569032442045
$ws.Range("N3").Value = 18.23483456406973

$ws.Range("B4").Value = 14.59109853562796
$ws.Range("C4").Value = 10.30221136720969
$ws.Range("D4").Value = 4.965917510826777
$ws.Range("F4").Value = 26.81820148456503
$ws.Range("G4").Value = 33.46245979942656
$ws.Range("H4").Value = 15.24354287142106
$ws.Range("I4").Value = 23.61347721000709
$ws.Range("L4").Value = 10.86202135188264
$ws.Range("M4").Value = 15.16731218757312
$ws.Range("N4").Value = 18.27717364278233

$ws.Range("B5").Value = 14.4675744116358
$ws.Range("C5").Value = 10.20858290751569
$ws.Range("D5").Value = 4.96509641906095
$ws.Range("F5").Value = 26.7982543615741
$ws.Range("G5").Value = 33.41277521257366
$ws.Range("H5").Value = 15.25197571479289
$ws.Range("I5").Value = 23.63516327514508
$ws.Range("L5").Value = 10.86692996637586
$ws.Range("M5").Value = 15.14825688947404
$ws.Range("N5").Value = 18.29489763923105

$ws.Range("B6").Value = 14.44700527082719
$ws.Range("C6").Value = 10.19291867488946
$ws.Range("D6").Value = 4.964961232461156
$ws.Range("F6").Value = 26.79507464553622
$ws.Range("G6").Value = 33.40471802920229
$ws.Range("H6").Value = 15.25342821153426
$ws.Range("I6").Value = 23.63885017062447
$ws.Range("L6").Value = 10.86776913822684
$ws.Range("M6").Value = 15.1451330219801
$ws.Range("N6").Value = 18.29786915551216

$ws.Range("B7").Value = 14.58943667318811
$ws.Range("C7").Value = 10.30095658002737
$ws.Range("D7").Value = 4.965906360078102
$ws.Range("F7").Value = 26.81792359696381
$ws.Range("G7").Value = 33.46177683168047
$ws.Range("H7").Value = 15.24365309690754
$ws.Range("I7").Value = 23.61376391122245
$ws.Range("L7").Value = 10.8620859355707
$ws.Range("M7").Value = 15.16705251226622
$ws.Range("N7").Value = 18.27741076786245

$ws.Range("B8").Value = 15.20809305775867
$ws.Range("C8").Value = 10.76007786051908
$ws.Range("D8").Value = 4.970242685496525
$ws.Range("F8").Value = 26.93791131313776
$ws.Range("G8").Value = 33.74127979083777
$ws.Range("H8").Value = 15.20812434685723
$ws.Range("I8").Value = 23.51559150469898
$ws.Range("L8").Value = 10.84050976494843
$ws.Range("M8").Value = 15.26952397660789
$ws.Range("N8").Value = 18.19125874561609

$ws.Range("B9").Value = 16.37347171382401
$ws.Range("C9").Value = 11.59223655185036
$ws.Range("D9").Value = 4.979321157225701
$ws.Range("F9").Value = 27.24089130024068
$ws.Range("G9").Value = 34.38537714873328
$ws.Range("H9").Value = 15.16591742959815
$ws.Range("I9").Value = 23.36843273800713
$ws.Range("L9").Value = 10.8107798187471
$ws.Range("M9").Value = 15.49005466915034
$ws.Range("N9").Value = 18.03704290897214

$ws.Range("B10").Value = 17.1908758352347
$ws.Range("C10").Value = 12.15846224946766
$ws.Range("D10").Value = 4.986324058140927
$ws.Range("F10").Value = 27.50341624393025
$ws.Range("G10").Value = 34.91330855392459
$ws.Range("H10").Value = 15.15174730016273
$ws.Range("I10").Value = 23.28825383375938
$ws.Range("L10").Value = 10.79659735028879
$ws.Range("M10").Value = 15.66299471769649
$ws.Range("N10").Value = 17.93261996911379

$ws.Range("B11").Value = 17.55256866878983
$ws.Range("C11").Value = 12.40562512250402
$ws.Range("D11").Value = 4.989578452455051
$ws.Range("F11").Value = 27.63116731956809
$ws.Range("G11").Value = 35.16443608934577
$ws.Range("H11").Value = 15.14897218648802
$ws.Range("I11").Value = 23.25790192863971
$ws.Range("L11").Value = 10.79180555523299
$ws.Range("M11").Value = 15.74380755154756
$ws.Range("N11").Value = 17.88702221980776

$ws.Range("B12").Value = 17.6879474961472
$ws.Range("C12").Value = 12.49767769076698
$ws.Range("D12").Value = 4.990820365954947
$ws.Range("F12").Value = 27.68070725545829
$ws.Range("G12").Value = 35.26102502619122
$ws.Range("H12").Value = 15.14844986295665
$ws.Range("I12").Value = 23.24729302372407
$ws.Range("L12").Value = 10.7902293262656
$ws.Range("M12").Value = 15.77469674464882
$ws.Range("N12").Value = 17.87002787612824

$ws.Range("B13").Value = 17.65886367728815
$ws.Range("C13").Value = 12.4779217870702
$ws.Range("D13").Value = 4.990552478552345
$ws.Range("F13").Value = 27.66998677135223
$ws.Range("G13").Value = 35.24015788487949
$ws.Range("H13").Value = 15.14853883986467
$ws.Range("I13").Value = 23.24953842705024
$ws.Range("L13").Value = 10.79055820262762
$ws.Range("M13").Value = 15.76803177052389
$ws.Range("N13").Value = 17.87367581272009

$ws.Range("B14").Value = 17.56373883948282
$ws.Range("C14").Value = 12.41322949646187
$ws.Range("D14").Value = 4.989680436077084
$ws.Range("F14").Value = 27.63521988813893
$ws.Range("G14").Value = 35.17235310838582
$ws.Range("H14").Value = 15.14891861893286
$ws.Range("I14").Value = 23.25701137279368
$ws.Range("L14").Value = 10.79167110399438
$ws.Range("M14").Value = 15.746343190486
$ws.Range("N14").Value = 17.8856186307261

$ws.Range("B15").Value = 17.50526199609681
$ws.Range("C15").Value = 12.37340146798354
$ws.Range("D15").Value = 4.989147517836753
$ws.Range("F15").Value = 27.61407462031816
$ws.Range("G15").Value = 35.13101249787173
$ws.Range("H15").Value = 15.14922009215347
$ws.Range("I15").Value = 23.26170410749558
$ws.Range("L15").Value = 10.79238381239852
$ws.Range("M15").Value = 15.7330950755573
$ws.Range("N15").Value = 17.8929694032114

$ws.Range("B16").Value = 17.16702296891215
$ws.Range("C16").Value = 12.14209601433011
$ws.Range("D16").Value = 4.986112730002554
$ws.Range("F16").Value = 27.49523212128922
$ws.Range("G16").Value = 34.89711032381386
$ws.Range("H16").Value = 15.15200259413012
$ws.Range("I16").Value = 23.29036097123453
$ws.Range("L16").Value = 10.79694387639152
$ws.Range("M16").Value = 15.65775473281232
$ws.Range("N16").Value = 17.93563810074843

$ws.Range("B17").Value = 16.95683021754385
$ws.Range("C17").Value = 11.99749608606909
$ws.Range("D17").Value = 4.984268350217596
$ws.Range("F17").Value = 27.4244354693027
$ws.Range("G17").Value = 34.75636906024148
$ws.Range("H17").Value = 15.15465036763205
$ws.Range("I17").Value = 23.30951219799558
$ws.Range("L17").Value = 10.80016622185868
$ws.Range("M17").Value = 15.61206970858768
$ws.Range("N17").Value = 17.9623008433738

$ws.Range("B18").Value = 16.83498438412391
$ws.Range("C18").Value = 11.91334792489677
$ws.Range("D18").Value = 4.983214002123826
$ws.Range("F18").Value = 27.38450110766103
$ws.Range("G18").Value = 34.67645618853174
$ws.Range("H18").Value = 15.15651880367099
$ws.Range("I18").Value = 23.32110345262678
$ws.Range("L18").Value = 10.80217589187105
$ws.Range("M18").Value = 15.58599574170109
$ws.Range("N18").Value = 17.97781591458396

$ws.Range("B19").Value = 16.79357066742994
$ws.Range("C19").Value = 11.8846902725368
$ws.Range("D19").Value = 4.98285814306247
$ws.Range("F19").Value = 27.37111601101071
$ws.Range("G19").Value = 34.64957983094676
$ws.Range("H19").Value = 15.15721073941997
$ws.Range("I19").Value = 23.32512683849269
$ws.Range("L19").Value = 10.80288317993124
$ws.Range("M19").Value = 15.57720301105652
$ws.Range("N19").Value = 17.9830998995106

$ws.Range("B20").Value = 16.97930472733818
$ws.Range("C20").Value = 12.01299056848937
$ws.Range("D20").Value = 4.984464018313046
$ws.Range("F20").Value = 27.4318907842946
$ws.Range("G20").Value = 34.7712444053715
$ws.Range("H20").Value = 15.15433274304306
$ws.Range("I20").Value = 23.30741387320608
$ws.Range("L20").Value = 10.79980702760123
$ws.Range("M20").Value = 15.61691211246976
$ws.Range("N20").Value = 17.95944399504615

$ws.Range("B21").Value = 17.59172330958003
$ws.Range("C21").Value = 12.43227338680802
$ws.Range("D21").Value = 4.989936319749264
$ws.Range("F21").Value = 27.64540047469296
$ws.Range("G21").Value = 35.19222920584782
$ws.Range("H21").Value = 15.1487927200869
$ws.Range("I21").Value = 23.25479234337838
$ws.Range("L21").Value = 10.79133775322722
$ws.Range("M21").Value = 15.75270602783371
$ws.Range("N21").Value = 17.88210334991791

$ws.Range("B22").Value = 17.98267803555742
$ws.Range("C22").Value = 12.69729028757458
$ws.Range("D22").Value = 4.993568262974951
$ws.Range("F22").Value = 27.7917049239385
$ws.Range("G22").Value = 35.47602420657
$ws.Range("H22").Value = 15.14825288441574
$ws.Range("I22").Value = 23.22555968605019
$ws.Range("L22").Value = 10.7871914844836
$ws.Range("M22").Value = 15.84311883859615
$ws.Range("N22").Value = 17.83314468065243

$ws.Range("B23").Value = 17.77490723291657
$ws.Range("C23").Value = 12.55668345295368
$ws.Range("D23").Value = 4.991624862565872
$ws.Range("F23").Value = 27.71301251310339
$ws.Range("G23").Value = 35.32379426777555
$ws.Range("H23").Value = 15.14825896375124
$ws.Range("I23").Value = 23.24068832357654
$ws.Range("L23").Value = 10.78927748172408
$ws.Range("M23").Value = 15.79471858853171
$ws.Range("N23").Value = 17.85912999856528

$ws.Range("B24").Value = 16.96914711644719
$ws.Range("C24").Value = 12.00598867343476
$ws.Range("D24").Value = 4.984375537998297
$ws.Range("F24").Value = 27.42851784281267
$ws.Range("G24").Value = 34.76451613682535
$ws.Range("H24").Value = 15.15447526274243
$ws.Range("I24").Value = 23.30836071645545
$ws.Range("L24").Value = 10.79996892994571
$ws.Range("M24").Value = 15.6147222652459
$ws.Range("N24").Value = 17.96073499570139

$ws.Range("B25").Value = 16.06437962617909
$ws.Range("C25").Value = 11.37481366097133
$ws.Range("D25").Value = 4.976805449631449
$ws.Range("F25").Value = 27.15180690144591
$ws.Range("G25").Value = 34.20122582373402
$ws.Range("H25").Value = 15.17438417475361
$ws.Range("I25").Value = 23.40335666808934
$ws.Range("L25").Value = 10.81747643131619
$ws.Range("M25").Value = 15.42840128529821
$ws.Range("N25").Value = 18.07719594383341
